$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "420.0"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "790.0"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4904.0"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "5751.0"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "7.24"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "4.390"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "12.148"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "14.064"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "0.100"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "2545.737"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "729.648"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "540.0"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1460.0"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2475.0"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "7303.0"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "9.459"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "8.299"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "3.525"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "11.387"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "0.034"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "0.151"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "2520.095"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "568.193"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0.32"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "19.79"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "89.0"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "15202.0"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1.28"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "220.31"
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "0.006"
$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "17080.89"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "0.36"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "0.3"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "25.66"
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "71.04"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "68.0"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "99.0"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "8579.0"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "13618.0"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "1.140"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "0.54"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "72.75"
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "114.445"
$ws.Range("O8").NumberFormat = "@"
$ws.Range("O8").Value = "34776.575"
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "13740.24"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "20.0"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "85.0"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "649.0"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "10242.0"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "0.31"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "10.3"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "162.57"
$ws.Range("O9").NumberFormat = "@"
$ws.Range("O9").Value = "54455.0"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "12049.41"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "34.0"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "89.0"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "2016.0"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "2229.0"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "0.77"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "0.67"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "9.416"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "10.434"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "0.053"
$ws.Range("O10").NumberFormat = "@"
$ws.Range("O10").Value = "11467.617"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "3632.269"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "114.0"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "0.32"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "0.05"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "3.4"
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = "1.86"
$ws.Range("M12").NumberFormat = "@"
$ws.Range("M12").Value = "0.031"
$ws.Range("N12").NumberFormat = "@"
$ws.Range("N12").Value = "0.061"
$ws.Range("O12").NumberFormat = "@"
$ws.Range("O12").Value = "1610.0"
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "1140.0"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "110.0"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "1235.0"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "1697.0"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "1.8"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "1.47"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = "20.24"
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L13").Value = "27.81"
$ws.Range("M13").NumberFormat = "@"
$ws.Range("M13").Value = "0.018"
$ws.Range("N13").NumberFormat = "@"
$ws.Range("N13").Value = "0.111"
$ws.Range("O13").NumberFormat = "@"
$ws.Range("O13").Value = "2665.45"
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "628.51"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "786.0"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "0.32"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "0.38"
$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = "0.7"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "12.88"
$ws.Range("N14").NumberFormat = "@"
$ws.Range("N14").Value = "0.093"
$ws.Range("O14").NumberFormat = "@"
$ws.Range("O14").Value = "4145.0"
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "1122.85"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "0.52"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "0.399"
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = "1.635"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "4.685"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "160.0"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "330.0"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "3336.0"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "2365.0"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "2.800"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "1.869"
$ws.Range("K16").NumberFormat = "@"
$ws.Range("K16").Value = "19.413"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "13.57"
$ws.Range("M16").NumberFormat = "@"
$ws.Range("M16").Value = "0.016"
$ws.Range("N16").NumberFormat = "@"
$ws.Range("N16").Value = "0.161"
$ws.Range("O16").NumberFormat = "@"
$ws.Range("O16").Value = "3199.81"
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "774.427"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "197.0"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "0.07"
$ws.Range("K17").NumberFormat = "@"
$ws.Range("K17").Value = "3.51"
$ws.Range("M17").NumberFormat = "@"
$ws.Range("M17").Value = "0.019"
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = "4925.0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "230.0"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "340.0"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "2689.0"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "2230.0"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "4.1"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "2.0"
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value = "12.027"
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value = "9.927"
$ws.Range("N18").NumberFormat = "@"
$ws.Range("N18").Value = "0.110"
$ws.Range("O18").NumberFormat = "@"
$ws.Range("O18").Value = "1970.172"
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = "643.255"
